# Update "想去人数" (F column) figures across all four sheets to match
# the refreshed data snapshot (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3306
$ws1.Range("F6").Value = 7726
$ws1.Range("F8").Value = 724
$ws1.Range("F9").Value = 1137
$ws1.Range("F11").Value = 1034
$ws1.Range("F12").Value = 11
$ws1.Range("F13").Value = 168
$ws1.Range("F14").Value = 1749
$ws1.Range("F15").Value = 363
$ws1.Range("F16").Value = 6179
$ws1.Range("F17").Value = 81
$ws1.Range("F18").Value = 2366
$ws1.Range("F19").Value = 106
$ws1.Range("F20").Value = 1026
$ws1.Range("F22").Value = 1034
$ws1.Range("F23").Value = 6310
$ws1.Range("F24").Value = 5962
$ws1.Range("F25").Value = 373
$ws1.Range("F26").Value = 165
$ws1.Range("F27").Value = 1078
$ws1.Range("F30").Value = 113
$ws1.Range("F32").Value = 1036
$ws1.Range("F33").Value = 108
$ws1.Range("F34").Value = 108
$ws1.Range("F35").Value = 85
$ws1.Range("F38").Value = 201
$ws1.Range("F39").Value = 79
$ws1.Range("F40").Value = 598
$ws1.Range("F42").Value = 338
$ws1.Range("F43").Value = 1222
$ws1.Range("F44").Value = 112
$ws1.Range("F45").Value = 477
$ws1.Range("F46").Value = 38
$ws1.Range("F47").Value = 3226
$ws1.Range("F48").Value = 100
$ws1.Range("F49").Value = 441
$ws1.Range("F50").Value = 51

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 81
$ws2.Range("F15").Value = 178
$ws2.Range("F24").Value = 6571
$ws2.Range("F35").Value = 32

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 1311
$ws3.Range("F8").Value = 549
$ws3.Range("F9").Value = 2138
$ws3.Range("F10").Value = 8904
$ws3.Range("F11").Value = 1002
$ws3.Range("F12").Value = 80

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 7726
$ws4.Range("F6").Value = 1311
$ws4.Range("F7").Value = 549
$ws4.Range("F8").Value = 2138
$ws4.Range("F9").Value = 1002
$ws4.Range("F11").Value = 80
$ws4.Range("F12").Value = 1137
$ws4.Range("F14").Value = 1034
$ws4.Range("F16").Value = 168
$ws4.Range("F18").Value = 363
$ws4.Range("F19").Value = 81
$ws4.Range("F20").Value = 2366
$ws4.Range("F21").Value = 1026
$ws4.Range("F22").Value = 1034
$ws4.Range("F23").Value = 6310
$ws4.Range("F24").Value = 5962
$ws4.Range("F25").Value = 373
$ws4.Range("F26").Value = 165
$ws4.Range("F27").Value = 1078
$ws4.Range("F30").Value = 113
$ws4.Range("F32").Value = 108
$ws4.Range("F33").Value = 108
$ws4.Range("F34").Value = 85
$ws4.Range("F36").Value = 201
$ws4.Range("F38").Value = 79
$ws4.Range("F39").Value = 598
$ws4.Range("F41").Value = 338
$ws4.Range("F44").Value = 112
$ws4.Range("F45").Value = 477
$ws4.Range("F46").Value = 3226
$ws4.Range("F47").Value = 100
$ws4.Range("F48").Value = 51
